$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.99999801046644
$ws.Range("E2").Value = 0.99999801046644

$ws.Range("D3").Value = 0.1818332035706449
$ws.Range("E3").Value = 0.1818332035706449

$ws.Range("D4").Value = [double]"1.912089694004327E-05"
$ws.Range("E4").Value = [double]"1.912089694004327E-05"

$ws.Range("D5").Value = [double]"3.415870495040572E-25"
$ws.Range("E5").Value = [double]"3.415870495040572E-25"

$ws.Range("D6").Value = [double]"7.037253649584015E-15"
$ws.Range("E6").Value = [double]"7.037253649584015E-15"

$ws.Range("D7").Value = 0.9999999999999494
$ws.Range("E7").Value = [double]"5.062616992290714E-14"

$ws.Range("D8").Value = [double]"7.021676074409973E-05"
$ws.Range("E8").Value = 0.9999297832392559

$ws.Range("D9").Value = 0.9976351319404879
$ws.Range("E9").Value = 0.002364868059512082

$ws.Range("D10").Value = [double]"2.529889376191617E-12"
$ws.Range("E10").Value = 0.9999999999974701

$ws.Range("D11").Value = [double]"2.218110811393918E-36"
$ws.Range("F11").Value = 13.16938495635986
